$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5458.978
$ws.Range("J17").Value = 5584.159
$ws.Range("L17").Value = 16752.477
$ws.Range("N17").Value = -17088.477
$ws.Range("H19").Value = 1468.2
$ws.Range("J19").Value = 1462
$ws.Range("L19").Value = 1462
$ws.Range("N19").Value = -1812
$ws.Range("H28").Value = 1892.6666
$ws.Range("I28").Value = 2030.4286
$ws.Range("K28").Value = 2030.4286
$ws.Range("M28").Value = -1545.4286
$ws.Range("H126").Value = 139999
$ws.Range("J126").Value = 139999
$ws.Range("L126").Value = 139999
$ws.Range("N126").Value = -149879
$ws.Range("H132").Value = 2488.5
$ws.Range("I132").Value = 2488.5
$ws.Range("K132").Value = 7465.5
$ws.Range("M132").Value = -4935.5
$ws.Range("H137").Value = 4339602.5
$ws.Range("J137").Value = 7940120.5
$ws.Range("L137").Value = 23820361.5
$ws.Range("N137").Value = -23825461.5
$ws.Range("H141").Value = 5654.483
$ws.Range("I141").Value = 4635.4546
$ws.Range("K141").Value = 13906.3638
$ws.Range("M141").Value = -8726.363799999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 173000
$ws.Range("J7").Value = 173000
$ws.Range("L7").Value = 173000
$ws.Range("N7").Value = -173228
$ws.Range("H44").Value = 69487
$ws.Range("J44").Value = 69974
$ws.Range("L44").Value = 69974
$ws.Range("N44").Value = -70950
$ws.Range("H52").Value = 59467
$ws.Range("I52").Value = 56986
$ws.Range("J52").Value = 61948
$ws.Range("K52").Value = 56986
$ws.Range("L52").Value = 61948
$ws.Range("M52").Value = -56668
$ws.Range("N52").Value = -62584
$ws.Range("H55").Value = 40013
$ws.Range("J55").Value = 69978
$ws.Range("L55").Value = 69978
$ws.Range("N55").Value = -70608
$ws.Range("H102").Value = 2334.6667
$ws.Range("I102").Value = 2101.6
$ws.Range("J102").Value = 3500
$ws.Range("K102").Value = 2101.6
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = -479.5999999999999
$ws.Range("N102").Value = -6744
$ws.Range("H132").Value = 419747.84
$ws.Range("I132").Value = 419747.84
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1259243.52
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1256713.52
$ws.Range("N132").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 306.73685
$ws.Range("I94").Value = 306.73685
$ws.Range("K94").Value = 306.73685
$ws.Range("M94").Value = 144.26315
$ws.Range("H99").Value = 2458.3845
$ws.Range("I99").Value = 1776.5
$ws.Range("K99").Value = 1776.5
$ws.Range("M99").Value = -278.5
$ws.Range("H105").Value = 3001.353
$ws.Range("I105").Value = 2503.1
$ws.Range("J105").Value = 3713.1428
$ws.Range("K105").Value = 2503.1
$ws.Range("L105").Value = 3713.1428
$ws.Range("M105").Value = -756.0999999999999
$ws.Range("N105").Value = -7207.1428
$ws.Range("H107").Value = 3600.375
$ws.Range("I107").Value = 3501.5264
$ws.Range("K107").Value = 3501.5264
$ws.Range("M107").Value = -1581.5264
$ws.Range("H134").Value = 7411573.5
$ws.Range("I134").Value = 7411573.5
$ws.Range("K134").Value = 22234720.5
$ws.Range("M134").Value = -22232185.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4457.5107
$ws.Range("I31").Value = 2024.6666
$ws.Range("J31").Value = 5597.9062
$ws.Range("K31").Value = 2024.6666
$ws.Range("L31").Value = 5597.9062
$ws.Range("M31").Value = -1729.6666
$ws.Range("N31").Value = -6187.9062
$ws.Range("H32").Value = 12000
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H34").Value = 4457.5107
$ws.Range("I34").Value = 2024.6666
$ws.Range("J34").Value = 5597.9062
$ws.Range("K34").Value = 2024.6666
$ws.Range("L34").Value = 5597.9062
$ws.Range("M34").Value = -1822.6666
$ws.Range("N34").Value = -6001.9062
$ws.Range("H69").Value = 39329
$ws.Range("I69").Value = 10091
$ws.Range("J69").Value = 47682.715
$ws.Range("K69").Value = 10091
$ws.Range("L69").Value = 47682.715
$ws.Range("M69").Value = -9342
$ws.Range("N69").Value = -49180.715
$ws.Range("H72").Value = 39329
$ws.Range("I72").Value = 10091
$ws.Range("J72").Value = 47682.715
$ws.Range("K72").Value = 30273
$ws.Range("L72").Value = 143048.145
$ws.Range("M72").Value = -26529
$ws.Range("N72").Value = -150536.145
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H132").Value = 1730.931
$ws.Range("I132").Value = 1806.7778
$ws.Range("J132").Value = 707
$ws.Range("K132").Value = 5420.3334
$ws.Range("L132").Value = 2121
$ws.Range("M132").Value = -2890.3334
$ws.Range("N132").Value = -7181
$ws.Range("H134").Value = 2669.2856
$ws.Range("I134").Value = 1982.6
$ws.Range("K134").Value = 5947.799999999999
$ws.Range("M134").Value = -3412.799999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 22.857143
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 60
$ws.Range("L2").Value = 600
$ws.Range("M2").Value = 53
$ws.Range("N2").Value = -826
$ws.Range("H11").Value = 500546.34
$ws.Range("I11").Value = 555839.9
$ws.Range("J11").Value = 334665.66
$ws.Range("K11").Value = 1667519.7
$ws.Range("L11").Value = 1003996.98
$ws.Range("M11").Value = -1667379.7
$ws.Range("N11").Value = -1004276.98
$ws.Range("H26").Value = 676
$ws.Range("I26").Value = 185
$ws.Range("J26").Value = 1658
$ws.Range("K26").Value = 555
$ws.Range("L26").Value = 4974
$ws.Range("M26").Value = -267
$ws.Range("N26").Value = -5550
$ws.Range("H136").Value = 1405.6
$ws.Range("I136").Value = 1405.6
$ws.Range("K136").Value = 4216.799999999999
$ws.Range("M136").Value = 883.2000000000007

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 3279.8
$ws.Range("I31").Value = 3349.75
$ws.Range("K31").Value = 3349.75
$ws.Range("M31").Value = -3057.75
$ws.Range("H37").Value = 3279.8
$ws.Range("I37").Value = 3349.75
$ws.Range("K37").Value = 3349.75
$ws.Range("M37").Value = -3072.75
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H51").Value = 107994.5
$ws.Range("J51").Value = 107994.5
$ws.Range("L51").Value = 107994.5
$ws.Range("N51").Value = -109012.5
$ws.Range("H97").Value = 1174.3
$ws.Range("I97").Value = 468
$ws.Range("K97").Value = 468
$ws.Range("M97").Value = 28
$ws.Range("H103").Value = 94998.25
$ws.Range("J103").Value = 94998.25
$ws.Range("L103").Value = 94998.25
$ws.Range("N103").Value = -97342.25
$ws.Range("H126").Value = 2945.5
$ws.Range("I126").Value = 2400
$ws.Range("J126").Value = 3491
$ws.Range("K126").Value = 7200
$ws.Range("L126").Value = 10473
$ws.Range("M126").Value = -4730
$ws.Range("N126").Value = -15413

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2301.3215
$ws.Range("J82").Value = 2158.5625
$ws.Range("L82").Value = 2158.5625
$ws.Range("N82").Value = -2880.5625
$ws.Range("H85").Value = 2301.3215
$ws.Range("J85").Value = 2158.5625
$ws.Range("L85").Value = 2158.5625
$ws.Range("N85").Value = -4654.5625
$ws.Range("H93").Value = 2580
$ws.Range("I93").Value = 2075
$ws.Range("J93").Value = 2782
$ws.Range("K93").Value = 2075
$ws.Range("L93").Value = 2782
$ws.Range("M93").Value = -827
$ws.Range("N93").Value = -5278
$ws.Range("H100").Value = 1487.5
$ws.Range("I100").Value = 1487.5
$ws.Range("K100").Value = 1487.5
$ws.Range("M100").Value = -946.5
$ws.Range("H132").Value = 1668334
$ws.Range("I132").Value = 5001002
$ws.Range("K132").Value = 15003006
$ws.Range("M132").Value = -15000476
$ws.Range("H136").Value = 19496
$ws.Range("I136").Value = 16744
$ws.Range("K136").Value = 50232
$ws.Range("M136").Value = -47682
$ws.Range("H139").Value = 76357.25
$ws.Range("J139").Value = 76357.25
$ws.Range("L139").Value = 76357.25
$ws.Range("N139").Value = -86637.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 333333400
$ws.Range("I14").Value = 333333400
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 333333400
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -333333232
$ws.Range("N14").ClearContents()
$ws.Range("H54").Value = 5444
$ws.Range("I54").Value = 5444
$ws.Range("K54").Value = 5444
$ws.Range("M54").Value = -4924
$ws.Range("H96").Value = 14994
$ws.Range("I96").Value = 13331.667
$ws.Range("J96").Value = 15825.167
$ws.Range("K96").Value = 13331.667
$ws.Range("L96").Value = 15825.167
$ws.Range("M96").Value = -11958.667
$ws.Range("N96").Value = -18571.167
$ws.Range("H100").Value = 1136.8235
$ws.Range("I100").Value = 786.6
$ws.Range("K100").Value = 1573.2
$ws.Range("M100").Value = -1032.2
$ws.Range("H126").Value = 3041.8
$ws.Range("I126").Value = 2176.625
$ws.Range("K126").Value = 6529.875
$ws.Range("M126").Value = -4059.875
$ws.Range("H136").Value = 2792.8
$ws.Range("I136").Value = 1854.9
$ws.Range("J136").Value = 4668.6
$ws.Range("K136").Value = 5564.700000000001
$ws.Range("L136").Value = 14005.8
$ws.Range("M136").Value = -3014.700000000001
$ws.Range("N136").Value = -19105.8
